# Atualização de bases das ligas, do dia: 2024-01-31 às 21-30
#
# The source data for several fixtures was shuffled: the match-detail
# columns (id, HomeTeam, AwayTeam, score, odds, ...) for a handful of rows
# were swapped/rotated among each other while the row's sequential index
# (column A) and everything else on the sheet stays put. Column A is the
# running id for the sheet and must NOT move with the rest of the row, so
# we only swap columns B:AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the B:AC content of every affected row BEFORE writing anything,
# so a row's new content never accidentally depends on another row that
# has already been overwritten in this same pass.
$rowNums = @(47, 48, 114, 115, 116, 117, 355, 357, 358, 600, 601, 837, 838, 839, 840)
$snapshot = @{}
foreach ($r in $rowNums) {
    $snapshot[$r] = $ws.Range("B" + $r + ":AC" + $r).Value2
}

# Permutation: each key row receives the snapshot taken from its value row.
$sourceFor = @{
    47  = 48
    48  = 47
    114 = 117
    115 = 116
    116 = 115
    117 = 114
    355 = 358
    357 = 355
    358 = 357
    600 = 601
    601 = 600
    837 = 839
    838 = 840
    839 = 837
    840 = 838
}

foreach ($r in $rowNums) {
    $src = $sourceFor[$r]
    $ws.Range("B" + $r + ":AC" + $r).Value = $snapshot[$src]
}
